# SPM_Updates_Performed.xlsx - log new entry:
#  - "fixed Railroad legend/fixed Area Office misspelling/added link to minute orders"
#    added as commit comment for the existing 11/20/2015 (row 7) entry
#  - a brand-new row 8 logging the Area Office renderer fix, dated 11/17/
#
# Note: values are written in the same order the cells appear in the
# original sharedStrings table (A8, then C7, then D8, then the already
# shared B8/F8 values) so new shared-string entries land in the expected
# order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 - date column (kept as literal text "11/17/", not a date serial)
$ws.Range("A8").Value = "11/17/"

# Row 7 - add the missing Commit Comment text
$ws.Range("C7").Value = "fixed Railroad legend/fixed Area Office misspelling/added link to minute orders"

# New row 8 - Description of Updates
$ws.Range("D8").Value = "Added renderer to Area Office layer to mimic transparency on old layer"

# New row 8 - Editor + QA/QC Testing Procedure (reuse existing values)
$ws.Range("B8").Value = "cbardash"
$ws.Range("F8").Value = "YES"

# Move the active selection to G7, matching the saved workbook state
$ws.Range("G7").Select()
